$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 : modality_1 now belongs to folder "13-protec-social"
$ws.Cells.Item(2,2).Value = "13-protec-social"

# Row 3 : a_or_b now belongs to folder "00-base"
$ws.Cells.Item(3,2).Value = "00-base"

# Row 4 : a_or_b_or_c now belongs to folder "00-base"
$ws.Cells.Item(4,2).Value = "00-base"

# Row 5 : only_zero now belongs to folder "00-base"
$ws.Cells.Item(5,2).Value = "00-base"

# New rows 6-9
$ws.Cells.Item(6,1).Value = "canton_sigle"
$ws.Cells.Item(6,2).Value = "00-base"
$ws.Cells.Item(6,3).Value = "Sigle des cantons"
$ws.Cells.Item(6,4).Value = "Liste des 26 cantons suisses et leur sigle"

$ws.Cells.Item(7,1).Value = "langue_sigle"
$ws.Cells.Item(7,2).Value = "00-base"
$ws.Cells.Item(7,3).Value = "Sigle des langues"
$ws.Cells.Item(7,4).Value = "Liste des 3 principales langues suisses et leur sigle"

$ws.Cells.Item(8,1).Value = "oui_non"
$ws.Cells.Item(8,2).Value = "00-base"
$ws.Cells.Item(8,3).Value = "Oui ou non"
$ws.Cells.Item(8,4).Value = "Oui ou non codé en 0 ou 1"

$ws.Cells.Item(9,1).Value = "vide"
$ws.Cells.Item(9,2).Value = "00-base"
$ws.Cells.Item(9,3).Value = "vide / manquant"
$ws.Cells.Item(9,4).Value = "Valeur vide ou manquante"

# Apply left/center alignment style to A6:A8 (matches new cellXfs index 1)
$seed = $ws.Cells.Item(6,1)
$seed.HorizontalAlignment = -4131
$seed.VerticalAlignment = -4108
$ws.Range("A6").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Resize the table to include the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D9"))

# Widen columns C and D to fit the new, longer content (matches Excel's
# post-edit auto-fit result: ~14.33 and 40 characters)
$ws.Columns.Item(3).ColumnWidth = 13.416666666666666
$ws.Columns.Item(4).ColumnWidth = 39.083333333333336

# Update selection to mirror the saved view state
$ws.Range("D10").Select()
